$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Date, Tested(all), Tested(daily), Positive(all),
# Positive(daily), Hospitalized, Intensive care, Discharged, Deaths(all), Deaths(daily))
$rows = @(
    @(76, 43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(77, 43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(78, 43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(79, 43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0)
)

# Use the last existing data row as the formatting template (it has the
# same banding/style pattern the new rows need) and copy it down first so
# the new rows inherit correct number formats / styles, then overwrite
# the values.
$template = $ws.Range("A75:J75")
foreach ($row in $rows) {
    $r = $row[0]
    $template.Copy($ws.Range("A" + $r + ":J" + $r))
}

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
    $ws.Cells.Item($r, 6).Value2 = $row[6]
    $ws.Cells.Item($r, 7).Value2 = $row[7]
    $ws.Cells.Item($r, 8).Value2 = $row[8]
    $ws.Cells.Item($r, 9).Value2 = $row[9]
    $ws.Cells.Item($r, 10).Value2 = $row[10]
}

# Extend the table ("Tabela1") so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J79"))

# Match the saved selection state from the diff.
[void]$ws.Range("A79:J79").Select()
